$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.270.50'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.594.40'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'211.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = "'0.501"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("D10").Value = "'19.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("D12").Value = '1.819.38'
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = '1.604.14'
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D15").Value = "'0.502"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.40%  '
$ws.Range("D16").Value = "'63.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '26.264.83'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = "'229.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'7.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.05%  '
$ws.Range("D20").Value = '0.0₃0719'
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").Value = "'4.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  +2.09%  '
$ws.Range("D24").Value = "'8.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = "'146.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.43%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = "'6.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("D29").Value = "'15.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '1.470.08'
$ws.Range("E32").Value = '  +4.22%  '
$ws.Range("E33").Value = '  +1.60%  '
$ws.Range("D34").Value = "'2.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("E36").Value = '  +0.72%  '
$ws.Range("E37").Value = '  -3.05%  '
$ws.Range("E38").Value = '  -0.60%  '
$ws.Range("D39").Value = "'0.817"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +1.27%  '
$ws.Range("D43").Value = "'0.927"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.90%  '
$ws.Range("D44").Value = '1.732.02'
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("D45").Value = "'0.753"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("D46").Value = "'60.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("D47").Value = "'88.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.26%  '
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("D50").Value = "'0.0953"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("E51").Value = '  -0.23%  '
